$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log row (row 11)
$ws.Range("A11").Value = "Hoi, hebben jullie al iets gehoord?"
$ws.Range("B11").Value = "mailmind.test@zohomail.eu"
$ws.Range("C11").Value = "Testmail #9: Hoi, hebben jullie al iets gehoord?"
$ws.Range("D11").Value = "Overig"
$ws.Range("E11").Value = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$ws.Range("F11").Value = "2025-07-31 21:41:15"
$ws.Range("G11").Value = "Ja"
$ws.Range("H11").Value = "Ja"
$ws.Range("I11").Value = "Nee"
$ws.Range("J11").Value = "Nee"

# Extend the conditional-formatting ranges to include the new row 11
$ws.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D11"))
$ws.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G11"))
$ws.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H11"))
$ws.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I11"))
$ws.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J11"))

# Update the Dashboard summary count for "Overig" (3 -> 4)
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 4
